$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 115.4136146666667
$ws.Range("H2").Value = 346.240844
$ws.Range("I2").Value = 0.2619217538490851
$ws.Range("J2").Value = 0.2619217538490851
$ws.Range("M2").Value = 77.08952333333333
$ws.Range("N2").Value = 231.26857
$ws.Range("O2").Value = 0.2403816673726824
$ws.Range("P2").Value = 0.2403816673726824
$ws.Range("Q2").Value = 8897.180540830341
$ws.Range("R2").Value = 80074.62486747307
$ws.Range("S2").Value = 0.06296118791142036
$ws.Range("T2").Value = 0.06296118791142037

# Row 3
$ws.Range("G3").Value = 115.4136146666667
$ws.Range("H3").Value = 346.240844
$ws.Range("I3").Value = 0.2619217538490851
$ws.Range("J3").Value = 0.2619217538490851
$ws.Range("O3").Value = 0.3167483425780597
$ws.Range("P3").Value = 0.3167483425780597
$ws.Range("Q3").Value = 11723.71928661495
$ws.Range("R3").Value = 105513.4735795345
$ws.Range("S3").Value = 0.08296328141683623
$ws.Range("T3").Value = 0.08296328141683623

# Row 4
$ws.Range("G4").Value = 115.4136146666667
$ws.Range("H4").Value = 346.240844
$ws.Range("I4").Value = 0.2619217538490851
$ws.Range("J4").Value = 0.2619217538490851
$ws.Range("O4").Value = 0.4428699900492579
$ws.Range("P4").Value = 0.4428699900492579
$ws.Range("Q4").Value = 16391.82513646118
$ws.Range("R4").Value = 147526.4262281506
$ws.Range("S4").Value = 0.1159972845208285
$ws.Range("T4").Value = 0.1159972845208285

# Row 5
$ws.Range("I5").Value = 0.6414314537852458
$ws.Range("J5").Value = 0.6414314537852458
$ws.Range("M5").Value = 77.08952333333333
$ws.Range("N5").Value = 231.26857
$ws.Range("O5").Value = 0.2403816673726824
$ws.Range("P5").Value = 0.2403816673726824
$ws.Range("Q5").Value = 21788.68828200824
$ws.Range("R5").Value = 196098.1945380741
$ws.Range("S5").Value = 0.1541883623661811
$ws.Range("T5").Value = 0.1541883623661811

# Row 6
$ws.Range("I6").Value = 0.6414314537852458
$ws.Range("J6").Value = 0.6414314537852458
$ws.Range("O6").Value = 0.3167483425780597
$ws.Range("P6").Value = 0.3167483425780597
$ws.Range("S6").Value = 0.2031723498639119
$ws.Range("T6").Value = 0.2031723498639119

# Row 7
$ws.Range("I7").Value = 0.6414314537852458
$ws.Range("J7").Value = 0.6414314537852458
$ws.Range("O7").Value = 0.4428699900492579
$ws.Range("P7").Value = 0.4428699900492579
$ws.Range("S7").Value = 0.2840707415551528
$ws.Range("T7").Value = 0.2840707415551528

# Row 8
$ws.Range("I8").Value = 0.09664679236566912
$ws.Range("J8").Value = 0.09664679236566913
$ws.Range("M8").Value = 77.08952333333333
$ws.Range("N8").Value = 231.26857
$ws.Range("O8").Value = 0.2403816673726824
$ws.Range("P8").Value = 0.2403816673726824
$ws.Range("Q8").Value = 3282.980308939716
$ws.Range("R8").Value = 29546.82278045745
$ws.Range("S8").Value = 0.02323211709508097
$ws.Range("T8").Value = 0.02323211709508098

# Row 9
$ws.Range("I9").Value = 0.09664679236566912
$ws.Range("J9").Value = 0.09664679236566913
$ws.Range("O9").Value = 0.3167483425780597
$ws.Range("P9").Value = 0.3167483425780597
$ws.Range("S9").Value = 0.03061271129731157
$ws.Range("T9").Value = 0.03061271129731157

# Row 10
$ws.Range("I10").Value = 0.09664679236566912
$ws.Range("J10").Value = 0.09664679236566913
$ws.Range("O10").Value = 0.4428699900492579
$ws.Range("P10").Value = 0.4428699900492579
$ws.Range("S10").Value = 0.04280196397327658
$ws.Range("T10").Value = 0.04280196397327658
